# Apply the "bienServiciosList" placeholder emphasis change:
#  - paragraph mark run-properties gain <w:i/><w:iCs/>
#  - the three runs making up "{bienServiciosList}" gain <w:i/><w:iCs/>
#    and their sz/szCs grow from 24 to 28 (half-points: 12pt -> 14pt)
$d = $word.ActiveDocument

# Locate the placeholder text; Find narrows $rng to the matched span.
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("{bienServiciosList}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find '{bienServiciosList}' in the document"
}

# Find the paragraph object that fully contains the matched range so we can
# rewrite it (InsertXML replaces the whole target range's contents, so we
# target the entire paragraph range - including its end-of-paragraph mark -
# to keep every original attribute/property intact).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -le $rng.Start -and $p.Range.End -ge $rng.End) {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    throw "Could not locate enclosing paragraph for '{bienServiciosList}'"
}
$prng = $target.Range

$newParaXml = '<w:p w14:paraId="36A68B45" w14:textId="35F87CB6" w:rsidR="00295D40" w:rsidRPr="005C6304" w:rsidRDefault="009C7F9B" w:rsidP="002A5396">' +
  '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:spacing w:after="160" w:line="259" w:lineRule="auto"/>' +
  '<w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:i/><w:iCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="es-NI"/></w:rPr>' +
  '</w:pPr>' +
  '<w:r w:rsidRPr="005C6304"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:i/><w:iCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="es-NI"/></w:rPr><w:t>{</w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r w:rsidRPr="005C6304"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:i/><w:iCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="es-NI"/></w:rPr><w:t>bienServiciosList</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r w:rsidRPr="005C6304"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:i/><w:iCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="es-NI"/></w:rPr><w:t>}</w:t></w:r>' +
  '</w:p>'

$xmlPkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
  '<w:body>' + $newParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$prng.InsertXML($xmlPkg)

Write-Host "Applied italic/size formatting to '{bienServiciosList}' placeholder."
